$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column-A cell style (bold, border, centered) to the newly added rows (31-37)
# before writing values, so the new rows match the formatting of existing player rows.
$ws.Range("A2").Copy($ws.Range("A31:A37"))

# Build the full data block for rows 1-37 (columns A:E) and write it in one shot
$data = New-Object "object[,]" 37,5
$data[0,0] = "player"
$data[0,1] = "W"
$data[0,2] = "L"
$data[0,3] = "RL"
$data[0,4] = "PTS"
$data[1,0] = "Rohan Chowla"
$data[1,1] = 19
$data[1,2] = 6
$data[1,3] = 1
$data[1,4] = 34
$data[2,0] = "Kevin Lee"
$data[2,1] = 16
$data[2,2] = 5
$data[2,3] = 0
$data[2,4] = 28
$data[3,0] = "Roman Ramirez"
$data[3,1] = 15
$data[3,2] = 5
$data[3,3] = 1
$data[3,4] = 22
$data[4,0] = "Jason Jackson"
$data[4,1] = 8
$data[4,2] = 7
$data[4,3] = 2
$data[4,4] = 15
$data[5,0] = "Aaron Carter"
$data[5,1] = 11
$data[5,2] = 7
$data[5,3] = 0
$data[5,4] = 14
$data[6,0] = "Kevin Cooper"
$data[6,1] = 10
$data[6,2] = 5
$data[6,3] = 0
$data[6,4] = 14
$data[7,0] = "Coby Lovelace"
$data[7,1] = 7
$data[7,2] = 5
$data[7,3] = 2
$data[7,4] = 13
$data[8,0] = "Nathan Snow"
$data[8,1] = 6
$data[8,2] = 3
$data[8,3] = 1
$data[8,4] = 11
$data[9,0] = "Cason Duszak"
$data[9,1] = 6
$data[9,2] = 5
$data[9,3] = 1
$data[9,4] = 11
$data[10,0] = "Gabe Silverstein"
$data[10,1] = 8
$data[10,2] = 8
$data[10,3] = 0
$data[10,4] = 10
$data[11,0] = "Will Simpson"
$data[11,1] = 7
$data[11,2] = 5
$data[11,3] = 1
$data[11,4] = 10
$data[12,0] = "Leah Baetcke"
$data[12,1] = 4
$data[12,2] = 6
$data[12,3] = 2
$data[12,4] = 10
$data[13,0] = "Jack Massingill"
$data[13,1] = 7
$data[13,2] = 8
$data[13,3] = 0
$data[13,4] = 9
$data[14,0] = "Eric LastName"
$data[14,1] = 3
$data[14,2] = 1
$data[14,3] = 1
$data[14,4] = 6
$data[15,0] = "Matthew Rusten"
$data[15,1] = 3
$data[15,2] = 3
$data[15,3] = 1
$data[15,4] = 6
$data[16,0] = "Ann Hall"
$data[16,1] = 4
$data[16,2] = 5
$data[16,3] = 1
$data[16,4] = 5
$data[17,0] = "Yvonne Nguyen"
$data[17,1] = 4
$data[17,2] = 4
$data[17,3] = 0
$data[17,4] = 5
$data[18,0] = "Carla Betancourt"
$data[18,1] = 3
$data[18,2] = 2
$data[18,3] = 0
$data[18,4] = 5
$data[19,0] = "Rose Roché"
$data[19,1] = 3
$data[19,2] = 6
$data[19,3] = 0
$data[19,4] = 4
$data[20,0] = "Luci Nguyen"
$data[20,1] = 2
$data[20,2] = 1
$data[20,3] = 1
$data[20,4] = 4
$data[21,0] = "Helen Dunn"
$data[21,1] = 1
$data[21,2] = 1
$data[21,3] = 1
$data[21,4] = 3
$data[22,0] = "Noah Dale"
$data[22,1] = 1
$data[22,2] = 3
$data[22,3] = 1
$data[22,4] = 3
$data[23,0] = "Reagan Fryatt"
$data[23,1] = 2
$data[23,2] = 2
$data[23,3] = 0
$data[23,4] = 2
$data[24,0] = "Kristian Banlaoi"
$data[24,1] = 1
$data[24,2] = 2
$data[24,3] = 0
$data[24,4] = 2
$data[25,0] = "Piper Parker"
$data[25,1] = 1
$data[25,2] = 4
$data[25,3] = 0
$data[25,4] = 2
$data[26,0] = "Anna Brown"
$data[26,1] = 1
$data[26,2] = 2
$data[26,3] = 0
$data[26,4] = 1
$data[27,0] = "Alex LastName"
$data[27,1] = 1
$data[27,2] = 2
$data[27,3] = 0
$data[27,4] = 1
$data[28,0] = "Julie Jackson"
$data[28,1] = 1
$data[28,2] = 2
$data[28,3] = 0
$data[28,4] = 1
$data[29,0] = "Carolyn LastName"
$data[29,1] = 1
$data[29,2] = 2
$data[29,3] = 0
$data[29,4] = 1
$data[30,0] = "Brian Tafazoli"
$data[30,1] = 0
$data[30,2] = 2
$data[30,3] = 0
$data[30,4] = 0
$data[31,0] = "Sam Tellis"
$data[31,1] = 0
$data[31,2] = 2
$data[31,3] = 0
$data[31,4] = 0
$data[32,0] = "Cassie Deering"
$data[32,1] = 0
$data[32,2] = 2
$data[32,3] = 0
$data[32,4] = 0
$data[33,0] = "Yafu LastName"
$data[33,1] = 0
$data[33,2] = 2
$data[33,3] = 0
$data[33,4] = 0
$data[34,0] = "Kim LastName"
$data[34,1] = 0
$data[34,2] = 2
$data[34,3] = 0
$data[34,4] = 0
$data[35,0] = "Evan Sooklal"
$data[35,1] = 0
$data[35,2] = 4
$data[35,3] = 0
$data[35,4] = 0
$data[36,0] = "Paul Bartenfeld"
$data[36,1] = 0
$data[36,2] = 8
$data[36,3] = 0
$data[36,4] = 0

$ws.Range("A1:E37").Value = $data
